$wb = $excel.ActiveWorkbook

# --- Rename worksheets (task order identifiers regenerated) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911110210803"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911150057838"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911150067856"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911150537822"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911151178281"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911109880383.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911110048633.csv"
$ws1.Range("B4").Value = "go_stims-16502911110058577.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911110200799.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16502911126407866.csv"
$ws2.Range("B3").Value = "OB-1650291112298637.csv"
$ws2.Range("B4").Value = "ZB-match_3-16502911117307928.csv"
$ws2.Range("B5").Value = "ZB-match_0-16502911111646838.csv"
$ws2.Range("B6").Value = "TB-16502911149827845.csv"
$ws2.Range("B7").Value = "TB-16502911147417824.csv"
$ws2.Range("B8").Value = "OB-16502911120066025.csv"
$ws2.Range("B9").Value = "ZB-match_9-16502911113343065.csv"
$ws2.Range("B10").Value = "OB-16502911123648126.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911150217834.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911150097883.csv"
$ws4.Range("B4").Value = "MM_stims-16502911150378237.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911150227823.csv"
$ws4.Range("B6").Value = "MM_stims-16502911150537822.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291115038788.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502911150857868.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911151017828.csv"
$ws5.Range("B4").Value = "SAT_stims-16502911150567834.csv"
$ws5.Range("B5").Value = "SAT_stims-16502911150697823.csv"
